$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Goal (per the diff): collapse the document back down to a single paragraph
# "This is lab 1 CS 109 Fall 2015" followed by the (now empty) _GoBack
# bookmark, dropping the two "Working on ..." list-bullet paragraphs and the
# blank paragraph between them, and dropping the List Paragraph style that
# only existed to support that bulleted list.
# ---------------------------------------------------------------------------

# 1) Insert a zero-width "_GoBack" bookmark right after the title text, while
#    it still precedes other content - this engine mis-places a bookmark
#    collapsed exactly on a paragraph-mark boundary, so we temporarily add a
#    placeholder character after the insertion point, bookmark the gap before
#    it, then remove the placeholder again.
$titleEnd = $d.Paragraphs.Item(1).Range.End - 1   # position right after "...2015"
$d.Range($titleEnd, $titleEnd).InsertAfter("X") | Out-Null
$gap = $d.Range($titleEnd, $titleEnd)
$d.Bookmarks.Add("_GoBack", $gap) | Out-Null
$d.Range($titleEnd, $titleEnd + 1).Delete() | Out-Null   # drop the "X" placeholder

# 2) Delete everything from the end of the title paragraph through the end of
#    the document body text (the blank paragraph + both bulleted paragraphs),
#    merging their trailing section break back onto the title paragraph.
$tailStart = $d.Paragraphs.Item(1).Range.End
$d.Range($tailStart, $d.Content.End).Delete() | Out-Null

# 3) The bulleted list paragraphs are gone, so the "List Paragraph" style that
#    only they used is now unused too - remove it, same as Word would offer
#    to clean up via "reduce file size" / unused style pruning.
foreach ($s in $d.Styles) {
    if ($s.NameLocal -eq "List Paragraph") {
        $s.Delete()
    }
}
